$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.653.96'
$ws.Range("E2").Value = '  +5.54%  '
$ws.Range("D3").Value = '3.108.58'
$ws.Range("E3").Value = '  +3.33%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '''584.74'
$ws.Range("E5").Value = '  +4.06%  '
$ws.Range("D6").Value = '''143.29'
$ws.Range("E6").Value = '  +3.46%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '3.099.45'
$ws.Range("E8").Value = '  +3.33%  '
$ws.Range("D9").Value = '''0.533'
$ws.Range("E9").Value = '  +1.68%  '
$ws.Range("E10").Value = '  +8.58%  '
$ws.Range("D11").Value = '''5.75'
$ws.Range("E11").Value = '  +9.61%  '
$ws.Range("E12").Value = '  +2.32%  '
$ws.Range("E13").Value = '  +5.06%  '
$ws.Range("D14").Value = '''35.57'
$ws.Range("E14").Value = '  +5.41%  '
$ws.Range("E15").Value = '  +0.71%  '
$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").Value = '''7.29'
$ws.Range("E16").Value = '  -0.45%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.105.78'
$ws.Range("E17").Value = '  +3.32%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '62.641.61'
$ws.Range("E18").Value = '  +5.32%  '
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").Value = '''454.15'
$ws.Range("E19").Value = '  +5.22%  '
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").Value = '''14.09'
$ws.Range("E20").Value = '  +2.53%  '
$ws.Range("B21").Value = 'Polygon'
$ws.Range("C21").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D21").Value = '''0.736'
$ws.Range("E21").Value = '  +1.82%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '''7.55'
$ws.Range("E22").Value = '  +5.70%  '
$ws.Range("B23").Value = 'InternetComputer(DFINITY)'
$ws.Range("C23").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D23").Value = '''13.70'
$ws.Range("E23").Value = '  +1.97%  '
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").Value = '''82.13'
$ws.Range("E24").Value = '  +1.65%  '
$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").Value = '''1.00'
$ws.Range("E25").Value = '  +0.01%  '
$ws.Range("B26").Value = 'ImmutableX'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D26").Value = '''2.27'
$ws.Range("E26").Value = '  +2.88%  '
$ws.Range("B27").Value = 'PancakeSwap'
$ws.Range("C27").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D27").Value = '''2.69'
$ws.Range("E27").Value = '  +5.47%  '
$ws.Range("B28").Value = 'FirstDigitalUSD'
$ws.Range("C28").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D28").Value = '''1.00'
$ws.Range("E28").Value = '  -0.27%  '
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").Value = '''8.26'
$ws.Range("E29").Value = '  +4.87%  '
$ws.Range("B30").Value = 'NEARProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D30").Value = '''6.87'
$ws.Range("E30").Value = '  +12.60%  '
$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").Value = '''0.112'
$ws.Range("E31").Value = '  +11.95%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").Value = '''27.09'
$ws.Range("E32").Value = '  +4.82%  '
$ws.Range("B33").Value = 'Mantle'
$ws.Range("C33").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D33").Value = '''1.05'
$ws.Range("E33").Value = '  +4.78%  '
$ws.Range("B34").Value = 'PEPE'
$ws.Range("C34").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D34").Value = '0.0₃0802'
$ws.Range("E34").Value = '  +4.88%  '
$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D35").Value = '''6.09'
$ws.Range("E35").Value = '  +1.78%  '
$ws.Range("B36").Value = 'Stacks'
$ws.Range("C36").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D36").Value = '''2.21'
$ws.Range("E36").Value = '  +4.56%  '
$ws.Range("B37").Value = 'OKB'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D37").Value = '''50.63'
$ws.Range("E37").Value = '  +3.60%  '
$ws.Range("B38").Value = 'dogwifhat'
$ws.Range("C38").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D38").Value = '''3.02'
$ws.Range("E38").Value = '  +9.82%  '
$ws.Range("B39").Value = 'Cosmos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D39").Value = '''8.81'
$ws.Range("E39").Value = '  +1.66%  '
$ws.Range("B40").Value = 'Bittensor'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D40").Value = '''424.28'
$ws.Range("E40").Value = '  +4.97%  '
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '2.954.13'
$ws.Range("E41").Value = '  +6.76%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").Value = '''0.0372'
$ws.Range("E42").Value = '  +5.15%  '
$ws.Range("B43").Value = 'TheGraph'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D43").Value = '''0.282'
$ws.Range("E43").Value = '  +11.43%  '
$ws.Range("B44").Value = 'Kaspa'
$ws.Range("C44").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D44").Value = '''0.111'
$ws.Range("E44").Value = '  +3.03%  '
$ws.Range("B45").Value = 'Fetch.AI'
$ws.Range("C45").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D45").Value = '''2.17'
$ws.Range("E45").Value = '  +7.52%  '
$ws.Range("B46").Value = 'Monero'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D46").Value = '''125.28'
$ws.Range("E46").Value = '  +1.67%  '
$ws.Range("B47").Value = 'USDe'
$ws.Range("C47").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D47").Value = '''0.999'
$ws.Range("E47").Value = '  -0.01%  '
$ws.Range("B48").Value = 'Arweave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D48").Value = '''35.04'
$ws.Range("E48").Value = '  -2.34%  '
$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").Value = '''0.112'
$ws.Range("E49").Value = '  +0.96%  '
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").Value = '''24.90'
$ws.Range("E50").Value = '  +5.65%  '
$ws.Range("B51").Value = 'ThetaToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D51").Value = '''2.18'
$ws.Range("E51").Value = '  +6.32%  '
